$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "cred"

# Column widths (closest values that round-trip to the target stored widths)
$ws.Columns.Item(1).ColumnWidth = 14.3
$ws.Columns.Item(2).ColumnWidth = 33.166666666666664

# Data (entered column-major to match the shared-string insertion order)
$ws.Range("A1").Value = "username"
$ws.Range("A2").Value = "password"
$ws.Range("B1").Value = "syprusgojek@gmail.com"
$ws.Range("B2").Value = "RoBotFra@432"

# Hyperlink on the password cell (mailto link, styled like a hyperlink)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:syprusgojek@gmail.com", "", "", "RoBotFra@432")

# Page setup (printed as A4/Letter by default; target uses paper size 9 = A4, portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
